$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "벌꿀 음료에 취한 레토는 비헌터에게 그녀가 예전에 겪었던 '썰'을 풀어주겠다고 한다.`n"
$ws.Range("D3").Value = "그중 어떤 일들은 학생자치단의 다른 멤버들도 잘 모르는 이야기였다.`n"
